$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New weekly column: AH, date 22_06_2021
$ws.Range("AH1").Value = "22_06_2021"

$ws.Range("AH2").Value = 231
$ws.Range("AH3").Value = 236
$ws.Range("AH4").Value = 772
$ws.Range("AH5").Value = 1166
$ws.Range("AH6").Value = 1698
$ws.Range("AH7").Value = 2542
$ws.Range("AH8").Value = 2437
$ws.Range("AH9").Value = 3265
$ws.Range("AH10").Value = 2577
$ws.Range("AH11").Value = 708

# Extend the weekly running-total sum into the new column
$ws.Range("AH12").Formula = "=SUM(AH2:AH11)"

# Update the view to match the new extent of data
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = $ws.Range("Y1").Column
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AG14").Select()
